$wb = $excel.ActiveWorkbook

# This script refreshes market-data-derived columns (H:N) across all eight
# job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to match the latest
# scheduled data-runner pull. Values only; no formulas/formatting involved.

# ----- ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 933.35297
$ws.Range("I6").Value = 446.625
$ws.Range("J6").Value = 1366
$ws.Range("K6").Value = 1339.875
$ws.Range("L6").Value = 4098
$ws.Range("M6").Value = -1227.875
$ws.Range("N6").Value = -4322

$ws.Range("H12").Value = 225
$ws.Range("I12").Value = 225
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 225
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -55
$ws.Range("N12").ClearContents()

$ws.Range("H87").Value = 15404.916
$ws.Range("J87").Value = 15404.916
$ws.Range("L87").Value = 15404.916
$ws.Range("N87").Value = -17900.916

$ws.Range("H90").Value = 15404.916
$ws.Range("J90").Value = 15404.916
$ws.Range("L90").Value = 46214.748
$ws.Range("N90").Value = -58694.748

$ws.Range("H129").Value = 728.5
$ws.Range("I129").Value = 529.875
$ws.Range("J129").Value = 993.3333
$ws.Range("K129").Value = 1589.625
$ws.Range("L129").Value = 2979.9999
$ws.Range("M129").Value = 3410.375
$ws.Range("N129").Value = -12979.9999

$ws.Range("H137").Value = 1453.1852
$ws.Range("I137").Value = 1588.3572
$ws.Range("K137").Value = 4765.071599999999
$ws.Range("M137").Value = -2215.071599999999

$ws.Range("H138").Value = 5615.1953
$ws.Range("J138").Value = 6021.3335
$ws.Range("L138").Value = 18064.0005
$ws.Range("N138").Value = -28344.0005

# ----- ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13391.929
$ws.Range("I32").Value = 10316.97
$ws.Range("K32").Value = 10316.97
$ws.Range("M32").Value = -10029.97

$ws.Range("H61").Value = 2124.9167
$ws.Range("I61").Value = 1873.174
$ws.Range("J61").Value = 2570.3076
$ws.Range("K61").Value = 1873.174
$ws.Range("L61").Value = 2570.3076
$ws.Range("M61").Value = -1661.174
$ws.Range("N61").Value = -2994.3076

$ws.Range("H136").Value = 2124.9167
$ws.Range("I136").Value = 1873.174
$ws.Range("J136").Value = 2570.3076
$ws.Range("K136").Value = 5619.522
$ws.Range("L136").Value = 7710.9228
$ws.Range("M136").Value = -3069.522
$ws.Range("N136").Value = -12810.9228

# ----- BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2843.7827
$ws.Range("I86").Value = 3100
$ws.Range("J86").Value = 2646.6924
$ws.Range("K86").Value = 3100
$ws.Range("L86").Value = 2646.6924
$ws.Range("M86").Value = -1977
$ws.Range("N86").Value = -4892.6924

$ws.Range("H89").Value = 2843.7827
$ws.Range("I89").Value = 3100
$ws.Range("J89").Value = 2646.6924
$ws.Range("K89").Value = 15500
$ws.Range("L89").Value = 13233.462
$ws.Range("M89").Value = -9884
$ws.Range("N89").Value = -24465.462

$ws.Range("H99").Value = 1985.8823
$ws.Range("I99").Value = 1184.75
$ws.Range("J99").Value = 2232.3845
$ws.Range("K99").Value = 1184.75
$ws.Range("L99").Value = 2232.3845
$ws.Range("M99").Value = 313.25
$ws.Range("N99").Value = -5228.3845

$ws.Range("H126").Value = 30966.666
$ws.Range("J126").Value = 30966.666
$ws.Range("L126").Value = 30966.666
$ws.Range("N126").Value = -40846.666

# ----- CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5538.643
$ws.Range("I31").Value = 2174.8572
$ws.Range("J31").Value = 12266.214
$ws.Range("K31").Value = 2174.8572
$ws.Range("L31").Value = 12266.214
$ws.Range("M31").Value = -1879.8572
$ws.Range("N31").Value = -12856.214

$ws.Range("H34").Value = 5538.643
$ws.Range("I34").Value = 2174.8572
$ws.Range("J34").Value = 12266.214
$ws.Range("K34").Value = 2174.8572
$ws.Range("L34").Value = 12266.214
$ws.Range("M34").Value = -1972.8572
$ws.Range("N34").Value = -12670.214

$ws.Range("H41").Value = 12591.8
$ws.Range("J41").Value = 30000
$ws.Range("L41").Value = 30000
$ws.Range("N41").Value = -30856

$ws.Range("H51").Value = 31068.562
$ws.Range("I51").Value = 10000
$ws.Range("J51").Value = 32473.133
$ws.Range("K51").Value = 10000
$ws.Range("L51").Value = 32473.133
$ws.Range("M51").Value = -9264
$ws.Range("N51").Value = -33945.133

$ws.Range("H58").Value = 2560.3333
$ws.Range("I58").Value = 1610.0555
$ws.Range("J58").Value = 3985.75
$ws.Range("K58").Value = 1610.0555
$ws.Range("L58").Value = 3985.75
$ws.Range("M58").Value = -1407.0555
$ws.Range("N58").Value = -4391.75

$ws.Range("H59").Value = 27800
$ws.Range("J59").Value = 33733.332
$ws.Range("L59").Value = 33733.332
$ws.Range("N59").Value = -36023.332

$ws.Range("H60").Value = 25950
$ws.Range("J60").Value = 27013.334
$ws.Range("L60").Value = 27013.334
$ws.Range("N60").Value = -28035.334

$ws.Range("H61").Value = 31068.562
$ws.Range("I61").Value = 10000
$ws.Range("J61").Value = 32473.133
$ws.Range("K61").Value = 10000
$ws.Range("L61").Value = 32473.133
$ws.Range("M61").Value = -9652
$ws.Range("N61").Value = -33169.133

$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498

$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H122").Value = 436394.7
$ws.Range("I122").Value = 1336.3334
$ws.Range("J122").Value = 911003.8
$ws.Range("K122").Value = 4009.0002
$ws.Range("L122").Value = 2733011.4
$ws.Range("M122").Value = -1559.0002
$ws.Range("N122").Value = -2737911.4

$ws.Range("H136").Value = 2560.3333
$ws.Range("I136").Value = 1610.0555
$ws.Range("J136").Value = 3985.75
$ws.Range("K136").Value = 4830.166499999999
$ws.Range("L136").Value = 11957.25
$ws.Range("M136").Value = -2280.166499999999
$ws.Range("N136").Value = -17057.25

# ----- CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 753.62964
$ws.Range("I34").Value = 300
$ws.Range("J34").Value = 771.0769
$ws.Range("K34").Value = 900
$ws.Range("L34").Value = 2313.2307
$ws.Range("M34").Value = -816
$ws.Range("N34").Value = -2481.2307

$ws.Range("H39").Value = 10074.786
$ws.Range("J39").Value = 3308.6667
$ws.Range("L39").Value = 9926.000100000001
$ws.Range("N39").Value = -10514.0001

$ws.Range("H55").Value = 4900.4443
$ws.Range("J55").Value = 5071.0586
$ws.Range("L55").Value = 15213.1758
$ws.Range("N55").Value = -15567.1758

$ws.Range("H113").Value = 1106.9048
$ws.Range("I113").Value = 1312.9286
$ws.Range("J113").Value = 694.8570999999999
$ws.Range("K113").Value = 3938.7858
$ws.Range("L113").Value = 2084.5713
$ws.Range("M113").Value = -1768.7858
$ws.Range("N113").Value = -6424.5713

$ws.Range("H122").Value = 992.7059
$ws.Range("I122").Value = 596.2308
$ws.Range("J122").Value = 2281.25
$ws.Range("K122").Value = 5366.077200000001
$ws.Range("L122").Value = 20531.25
$ws.Range("M122").Value = -2916.077200000001
$ws.Range("N122").Value = -25431.25

$ws.Range("H131").Value = 978.74
$ws.Range("J131").Value = 1001.5158
$ws.Range("L131").Value = 3004.5474
$ws.Range("N131").Value = -13084.5474

$ws.Range("H132").Value = 2352.5
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2352.5
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 21172.5
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -26232.5

# ----- GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2235
$ws.Range("I80").Value = 1902.5
$ws.Range("K80").Value = 1902.5
$ws.Range("M80").Value = -904.5

$ws.Range("H83").Value = 2235
$ws.Range("I83").Value = 1902.5
$ws.Range("K83").Value = 9512.5
$ws.Range("M83").Value = -4520.5

$ws.Range("H122").Value = 2365.0571
$ws.Range("I122").Value = 1704.8182
$ws.Range("K122").Value = 5114.4546
$ws.Range("M122").Value = -2664.4546

$ws.Range("H126").Value = 2171.8518
$ws.Range("I126").Value = 1944.9231
$ws.Range("K126").Value = 5834.7693
$ws.Range("M126").Value = -3364.7693

# ----- LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 19611172
$ws.Range("I136").Value = 3738.6924
$ws.Range("J136").Value = 83335336
$ws.Range("K136").Value = 11216.0772
$ws.Range("L136").Value = 250006008
$ws.Range("M136").Value = -8666.0772
$ws.Range("N136").Value = -250011108

# ----- WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2416.6667
$ws.Range("I81").Value = 2400
$ws.Range("J81").Value = 2428.5715
$ws.Range("K81").Value = 4800
$ws.Range("L81").Value = 4857.143
$ws.Range("M81").Value = -3739
$ws.Range("N81").Value = -6979.143

$ws.Range("H84").Value = 2416.6667
$ws.Range("I84").Value = 2400
$ws.Range("J84").Value = 2428.5715
$ws.Range("K84").Value = 24000
$ws.Range("L84").Value = 24285.715
$ws.Range("M84").Value = -18696
$ws.Range("N84").Value = -34893.715

$ws.Range("H126").Value = 111333.9
$ws.Range("I126").Value = 158004.86
$ws.Range("J126").Value = 2435
$ws.Range("K126").Value = 474014.58
$ws.Range("L126").Value = 7305
$ws.Range("M126").Value = -471544.58
$ws.Range("N126").Value = -12245

$ws.Range("H136").Value = 6064.86
$ws.Range("I136").Value = 5057
$ws.Range("K136").Value = 15171
$ws.Range("M136").Value = -12621
